$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Clear the "-" placeholder text in C17, C31, C34 (now blank cells)
$ws.Range("C17").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("C34").Value = ""

# Update revised monthly figures
$ws.Range("AV11").Value = 228
$ws.Range("AX11").Value = 92
$ws.Range("BA11").Value = 163
$ws.Range("AV12").Value = 485
$ws.Range("AY12").Value = 353
$ws.Range("BA12").Value = 343
$ws.Range("AX13").Value = 31
$ws.Range("BA13").Value = 19
$ws.Range("AV14").Value = 126
$ws.Range("BA14").Value = 120
$ws.Range("AV18").Value = 868
$ws.Range("AX18").Value = 709
$ws.Range("AY18").Value = 585
$ws.Range("BA18").Value = 645
$ws.Range("AW25").Value = 242
$ws.Range("AX25").Value = 212
$ws.Range("AY25").Value = 41
$ws.Range("AY26").Value = 395
$ws.Range("AY27").Value = 24
$ws.Range("AY28").Value = 114
$ws.Range("AW35").Value = 876
$ws.Range("AX35").Value = 918
$ws.Range("AY35").Value = 574
$ws.Range("AW42").Value = 38106528
$ws.Range("AX42").Value = 33644300
$ws.Range("AY42").Value = 8156854
$ws.Range("AX43").Value = 113965452
$ws.Range("AY43").Value = 79887387
$ws.Range("AY44").Value = 7601381
$ws.Range("AY45").Value = 26010864
$ws.Range("AW54").Value = 173500478
$ws.Range("AY54").Value = 121656486
$ws.Range("AW61").Value = 157424783326
$ws.Range("AX61").Value = 159147876817
$ws.Range("AY61").Value = 180844382905
$ws.Range("AX62").Value = 210268360500
$ws.Range("AY62").Value = 201455690433
$ws.Range("AY63").Value = 306628671711
$ws.Range("AY64").Value = 231300256339
